$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Pris Karla" -> "Karla Pris"
$ws.Range("B1").Value = "Karla Pris"

# Add new header "Viggo Pris" for column D (was duplicate "Pris Karla")
$ws.Range("D1").Value = "Viggo Pris"

# Update the active cell selection on the sheet
$ws.Range("H8").Select()
